$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.563.75'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.623.84'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.527'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0611'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0889'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.854.36'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.628.48'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.529.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.74'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('E22').Value = '  +2.77%  '
$ws.Range('E24').Value = '  +5.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.465.73'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E34').Value = '  -2.52%  '
$ws.Range('E35').Value = '  -2.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.940'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.875'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.66%  '
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.556'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.68%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '67.53'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('B44').Value = 'mCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.05%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.76'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.01%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.764.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.48'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0105'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.85%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.100'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.48%  '
